$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on all target Price (D) cells so that
# numeric-looking strings (e.g. "74.00", "0.0000100") are written as
# literal text, matching the original inlineStr cell contents, instead
# of being coerced into numeric values by Excel.
$priceCells = @(
    'D2', 'D3', 'D5', 'D7', 'D9', 'D10', 'D11', 'D12', 'D13', 'D15', 'D16', 'D18', 'D19', 'D20', 'D22', 'D23', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D33', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D44', 'D45', 'D46', 'D51'
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the cell value updates from the diff
$ws.Range('D2').Value = '42.618.10'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '2.282.26'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '251.64'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  +2.56%  '
$ws.Range('D7').Value = '74.00'
$ws.Range('E7').Value = '  +9.61%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '0.639'
$ws.Range('E9').Value = '  +3.81%  '
$ws.Range('D10').Value = '39.36'
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('D11').Value = '0.0986'
$ws.Range('E11').Value = '  +5.46%  '
$ws.Range('D12').Value = '59.01'
$ws.Range('E12').Value = '  -0.49%  '
$ws.Range('D13').Value = '7.34'
$ws.Range('E13').Value = '  +4.46%  '
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').Value = '2.624.34'
$ws.Range('E15').Value = '  +3.46%  '
$ws.Range('D16').Value = '15.08'
$ws.Range('E16').Value = '  +4.17%  '
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '2.279.92'
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('D19').Value = '42.553.36'
$ws.Range('E19').Value = '  +1.84%  '
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +4.86%  '
$ws.Range('E21').Value = '  +2.99%  '
$ws.Range('D22').Value = '72.34'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '232.65'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('E24').Value = '  +10.46%  '
$ws.Range('D25').Value = '3.92'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = '11.49'
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '2.41'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('D29').Value = '3.62'
$ws.Range('E29').Value = '  -1.62%  '
$ws.Range('D30').Value = '2.14'
$ws.Range('E30').Value = '  +3.38%  '
$ws.Range('D31').Value = '166.94'
$ws.Range('E31').Value = '  +0.42%  '
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').Value = '6.38'
$ws.Range('E33').Value = '  +8.79%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '32.46'
$ws.Range('E35').Value = '  +25.97%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.0817'
$ws.Range('E36').Value = '  +5.23%  '
$ws.Range('D37').Value = '0.125'
$ws.Range('E37').Value = '  +2.81%  '
$ws.Range('D38').Value = '4.72'
$ws.Range('E38').Value = '  +16.47%  '
$ws.Range('D39').Value = '4.74'
$ws.Range('E39').Value = '  +3.68%  '
$ws.Range('D40').Value = '0.0306'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').Value = '13.89'
$ws.Range('E41').Value = '  +15.97%  '
$ws.Range('E42').Value = '  +5.66%  '
$ws.Range('E43').Value = '  +6.26%  '
$ws.Range('D44').Value = '0.213'
$ws.Range('E44').Value = '  +8.92%  '
$ws.Range('D45').Value = '9.14'
$ws.Range('E45').Value = '  +6.71%  '
$ws.Range('D46').Value = '61.98'
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('E47').Value = '  -7.17%  '
$ws.Range('E48').Value = '  +4.17%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('D51').Value = '97.85'
$ws.Range('E51').Value = '  +5.63%  '

# Reset style on the Price cells back to Normal so no stray
# cell-format (style index) is left attached to these cells -
# matches the original (unstyled) inline-string cells.
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
